$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.032735120064274
$ws.Range("D2").Value = 0.2833460919314064
$ws.Range("E2").Value = 0.1393265427079342
$ws.Range("F2").Value = 6.59407392202786
$ws.Range("G2").Value = 0.002685557512994389
$ws.Range("J2").Value = 0.1342384014459554
$ws.Range("L2").Value = 1.545165939102986
$ws.Range("B3").Value = 2.955087695476834
$ws.Range("D3").Value = 0.2538119477628413
$ws.Range("E3").Value = 0.1208143278345304
$ws.Range("F3").Value = 6.485554242696054
$ws.Range("G3").Value = 0.002696892544302101
$ws.Range("J3").Value = 0.1185460403325322
$ws.Range("L3").Value = 1.47854960731118
$ws.Range("B4").Value = 2.909412039891492
$ws.Range("D4").Value = 0.2358781291054299
$ws.Range("E4").Value = 0.1094861851150384
$ws.Range("F4").Value = 6.423296494313178
$ws.Range("G4").Value = 0.002704200464143268
$ws.Range("J4").Value = 0.1088865404359893
$ws.Range("L4").Value = 1.438808411844349
$ws.Range("B5").Value = 2.891299110976604
$ws.Range("D5").Value = 0.2286177645507337
$ws.Range("E5").Value = 0.1048784935284743
$ws.Range("F5").Value = 6.399008706135959
$ws.Range("G5").Value = 0.002707266470540425
$ws.Range("J5").Value = 0.1049435822378371
$ws.Range("L5").Value = 1.422902673287439
$ws.Range("B6").Value = 2.888321612261223
$ws.Range("D6").Value = 0.2274150007617095
$ws.Range("E6").Value = 0.1041138803811137
$ws.Range("F6").Value = 6.395040639060682
$ws.Range("G6").Value = 0.002707780904018687
$ws.Range("J6").Value = 0.1042884417644956
$ws.Range("L6").Value = 1.42027891608214
$ws.Range("B7").Value = 2.909165740455705
$ws.Range("D7").Value = 0.2357800230307134
$ws.Range("E7").Value = 0.1094240106403674
$ws.Range("F7").Value = 6.422964579028985
$ws.Range("G7").Value = 0.002704241456565392
$ws.Range("J7").Value = 0.1088333918475968
$ws.Range("L7").Value = 1.438592734645795
$ws.Range("B8").Value = 3.005545458235076
$ws.Range("D8").Value = 0.2731195988216939
$ws.Range("E8").Value = 0.1329349113310556
$ws.Range("F8").Value = 6.555738157996103
$ws.Range("G8").Value = 0.002689393824784348
$ws.Range("J8").Value = 0.1288324017015015
$ws.Range("L8").Value = 1.521953720532196
$ws.Range("B9").Value = 3.210558293627059
$ws.Range("D9").Value = 0.348051583508493
$ws.Range("E9").Value = 0.1793951981833999
$ws.Range("F9").Value = 6.851610330379515
$ws.Range("G9").Value = 0.002663020886258381
$ws.Range("J9").Value = 0.1678836278594389
$ws.Range("L9").Value = 1.694789860026503
$ws.Range("B10").Value = 3.37116991353588
$ws.Range("D10").Value = 0.4043204471635136
$ws.Range("E10").Value = 0.2138212858788791
$ws.Range("F10").Value = 7.091783433543412
$ws.Range("G10").Value = 0.00264528981801868
$ws.Range("J10").Value = 0.1965136111849546
$ws.Range("L10").Value = 1.827713469051957
$ws.Range("B11").Value = 3.446458905480597
$ws.Range("D11").Value = 0.4302218893978136
$ws.Range("E11").Value = 0.2295630426927886
$ws.Range("F11").Value = 7.206251961215912
$ws.Range("G11").Value = 0.002637574814261337
$ws.Range("J11").Value = 0.2095344441454614
$ws.Range("L11").Value = 1.889527930434724
$ws.Range("B12").Value = 3.475293200260182
$ws.Range("D12").Value = 0.4400770719034313
$ws.Range("E12").Value = 0.235537095412667
$ws.Range("F12").Value = 7.250368939011253
$ws.Range("G12").Value = 0.00263470334301912
$ws.Range("J12").Value = 0.2144654199861691
$ws.Range("L12").Value = 1.913133479371822
$ws.Range("B13").Value = 3.469068745300206
$ws.Range("D13").Value = 0.4379524474188088
$ws.Range("E13").Value = 0.2342498769466772
$ws.Range("F13").Value = 7.24083296745863
$ws.Range("G13").Value = 0.002635319547450723
$ws.Range("J13").Value = 0.2134034223934833
$ws.Range("L13").Value = 1.908040739849298
$ws.Range("B14").Value = 3.44882460403619
$ws.Range("D14").Value = 0.4310317238721666
$ws.Range("E14").Value = 0.2300542635066449
$ws.Range("F14").Value = 7.209865934689958
$ws.Range("G14").Value = 0.002637337576371402
$ws.Range("J14").Value = 0.2099401088843962
$ws.Range("L14").Value = 1.891465985809646
$ws.Range("B15").Value = 3.436466794839987
$ws.Range("D15").Value = 0.4267987776844677
$ws.Range("E15").Value = 0.2274860609661999
$ws.Range("F15").Value = 7.190998662217112
$ws.Range("G15").Value = 0.002638580180273059
$ws.Range("J15").Value = 0.2078187861961709
$ws.Range("L15").Value = 1.88133934817148
$ws.Range("B16").Value = 3.366294721375766
$ws.Range("D16").Value = 0.4026341154180386
$ws.Range("E16").Value = 0.2127942699596161
$ws.Range("F16").Value = 7.084409291194135
$ws.Range("G16").Value = 0.002645801036520462
$ws.Range("J16").Value = 0.1956626670866655
$ws.Range("L16").Value = 1.823701194883142
$ws.Range("B17").Value = 3.323819258007177
$ws.Range("D17").Value = 0.3878900259778675
$ws.Range("E17").Value = 0.2038030431400841
$ws.Range("F17").Value = 7.020369770645374
$ws.Range("G17").Value = 0.002650320377536516
$ws.Range("J17").Value = 0.1882049862614537
$ws.Range("L17").Value = 1.788689855527139
$ws.Range("B18").Value = 3.299597767465684
$ws.Range("D18").Value = 0.3794380824627694
$ws.Range("E18").Value = 0.1986390843023003
$ws.Range("F18").Value = 6.984025297289406
$ws.Range("G18").Value = 0.002652952841984115
$ws.Range("J18").Value = 0.1839151897449653
$ws.Range("L18").Value = 1.768678774853527
$ws.Range("B19").Value = 3.291432614266
$ws.Range("D19").Value = 0.3765811961135057
$ws.Range("E19").Value = 0.1968919150204016
$ws.Range("F19").Value = 6.971803110831218
$ws.Range("G19").Value = 0.002653849840076113
$ws.Range("J19").Value = 0.1824626595521721
$ws.Range("L19").Value = 1.76192497482441
$ws.Range("B20").Value = 3.328319157099486
$ws.Range("D20").Value = 0.3894565911713812
$ws.Range("E20").Value = 0.2047593830444043
$ws.Range("F20").Value = 7.02713609148185
$ws.Range("G20").Value = 0.002649835867735606
$ws.Range("J20").Value = 0.1889988988156546
$ws.Range("L20").Value = 1.792403749544007
$ws.Range("B21").Value = 3.454761974041219
$ws.Range("D21").Value = 0.4330632110989541
$ws.Range("E21").Value = 0.2312862535425637
$ws.Range("F21").Value = 7.21894063164558
$ws.Range("G21").Value = 0.002636743477980861
$ws.Range("J21").Value = 0.2109573548461583
$ws.Range("L21").Value = 1.89632899238245
$ws.Range("B22").Value = 3.53929006693096
$ws.Range("D22").Value = 0.4618373397009918
$ws.Range("E22").Value = 0.2486994660597333
$ws.Range("F22").Value = 7.348795457957067
$ws.Range("G22").Value = 0.002628478271911649
$ws.Range("J22").Value = 0.2253104008400868
$ws.Range("L22").Value = 1.965404841225791
$ws.Range("B23").Value = 3.49400151988732
$ws.Range("D23").Value = 0.4464538689914832
$ws.Range("E23").Value = 0.2393982701261592
$ws.Range("F23").Value = 7.279070822161088
$ws.Range("G23").Value = 0.002632863046298077
$ws.Range("J23").Value = 0.2176494884764679
$ws.Range("L23").Value = 1.928430707119105
$ws.Range("B24").Value = 3.32628413560326
$ws.Range("D24").Value = 0.3887482705630418
$ws.Range("E24").Value = 0.2043270056224316
$ws.Range("F24").Value = 7.024075567033719
$ws.Range("G24").Value = 0.002650054807790383
$ws.Range("J24").Value = 0.1886399780942014
$ws.Range("L24").Value = 1.790724332675438
$ws.Range("B25").Value = 3.153359836800007
$ws.Range("D25").Value = 0.3275797610137658
$ws.Range("E25").Value = 0.1667808001507325
$ws.Range("F25").Value = 6.767650721129201
$ws.Range("G25").Value = 0.002669864586587028
$ws.Range("J25").Value = 0.1573330014196586
$ws.Range("L25").Value = 1.647007169176959
